$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicated rows 8-13 (the data was repeated; only the first
# occurrence, rows 2-7, is kept).
$ws.Rows("8:13").Delete()

# Refresh the uuid value in column G (uuid) for the remaining rows.
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 7).Value = "44bb2153-845f-4543-8e0b-e127667e7e30"
}

# Apply an integer number format to the downtime column (H) for the
# remaining data rows, as used by the new highcharts visualization.
$ws.Range("H2:H7").NumberFormat = "0"
